# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0254748296512313
$ws.Range("C2").Value = 1.487060150505943
$ws.Range("D2").Value = 8.204358864081295
$ws.Range("E2").Value = 2.86432520222151
$ws.Range("F2").Value = 2.92858426574186
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.6231471476049726
$ws.Range("C3").Value = 2.149823722789349
$ws.Range("D3").Value = 19.09538388526634
$ws.Range("E3").Value = 4.369826528051926
$ws.Range("F3").Value = 4.426949790444016
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.8905896291048404
$ws.Range("C4").Value = 1.378004675742914
$ws.Range("D4").Value = 7.676892968144138
$ws.Range("E4").Value = 2.770720658627307
$ws.Range("F4").Value = 2.688481027403539
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.09339861389089041
$ws.Range("C5").Value = 0.4914936818093819
$ws.Range("D5").Value = 0.4832682964296551
$ws.Range("E5").Value = 0.6951750113673931
$ws.Range("F5").Value = 0.7067680296532315
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.1408622386988335
$ws.Range("C6").Value = 0.8298683219479311
$ws.Range("D6").Value = 1.789059693641156
$ws.Range("E6").Value = 1.337557360878836
$ws.Range("F6").Value = 1.366567739176598
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = -0.06671830263824358
$ws.Range("C7").Value = 0.6750367220520066
$ws.Range("D7").Value = 1.315810782452338
$ws.Range("E7").Value = 1.147087957591892
$ws.Range("F7").Value = 1.178345552778965
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = -0.01304991229360516
$ws.Range("C8").Value = 0.5207978849898635
$ws.Range("D8").Value = 0.4860254750586991
$ws.Range("E8").Value = 0.6971552732775527
$ws.Range("F8").Value = 0.7184852978842488
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = 0.276249544383158
$ws.Range("C9").Value = 0.5074605133255659
$ws.Range("D9").Value = 0.4401873792382052
$ws.Range("E9").Value = 0.6634661854519831
$ws.Range("F9").Value = 0.6230022522945293
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.1691041117434205
$ws.Range("C10").Value = 0.3720285681882072
$ws.Range("D10").Value = 0.2384939909285282
$ws.Range("E10").Value = 0.4883584656054692
$ws.Range("F10").Value = 0.4742262009089842
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = 0.2048918841260996
$ws.Range("C11").Value = 0.4286216204399972
$ws.Range("D11").Value = 0.2971361016417544
$ws.Range("E11").Value = 0.545101918581979
$ws.Range("F11").Value = 0.5257550464960803
